$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.577.54"
$ws.Range("E2").Value = "  +2.37%  "
$ws.Range("D3").Value = "1.996.08"
$ws.Range("E3").Value = "  +6.24%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("D4").Style = $ws.Range("C4").Style
$ws.Range("E4").Value = "  -0.15%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "327.64"
$ws.Range("D5").Style = $ws.Range("C5").Style
$ws.Range("E5").Value = "  +1.01%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9998"
$ws.Range("D6").Style = $ws.Range("C6").Style
$ws.Range("E6").Value = "  -0.13%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4684"
$ws.Range("D7").Style = $ws.Range("C7").Style
$ws.Range("E7").Value = "  +1.46%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3948"
$ws.Range("D8").Style = $ws.Range("C8").Style
$ws.Range("E8").Value = "  +1.77%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.08033"
$ws.Range("D9").Style = $ws.Range("C9").Style
$ws.Range("E9").Value = "  +2.26%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.001"
$ws.Range("D10").Style = $ws.Range("C10").Style
$ws.Range("E10").Value = "  +1.80%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "22.89"
$ws.Range("D11").Style = $ws.Range("C11").Style
$ws.Range("E11").Value = "  +5.12%  "
$ws.Range("D12").Value = "1.988.32"
$ws.Range("E12").Value = "  +6.09%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "7.241"
$ws.Range("D13").Style = $ws.Range("C13").Style
$ws.Range("E13").Value = "  +3.59%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.870"
$ws.Range("D14").Style = $ws.Range("C14").Style
$ws.Range("E14").Value = "  +3.58%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.07132"
$ws.Range("D15").Style = $ws.Range("C15").Style
$ws.Range("E15").Value = "  +2.16%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "88.98"
$ws.Range("D16").Style = $ws.Range("C16").Style
$ws.Range("E16").Value = "  +0.54%  "
$ws.Range("E17").Value = "  -0.09%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001004"
$ws.Range("D18").Style = $ws.Range("C18").Style
$ws.Range("E18").Value = "  +1.02%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "17.42"
$ws.Range("D19").Style = $ws.Range("C19").Style
$ws.Range("E19").Value = "  +2.73%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.000"
$ws.Range("D20").Style = $ws.Range("C20").Style
$ws.Range("E20").Value = "  -0.11%  "
$ws.Range("D21").Value = "29.586.90"
$ws.Range("E21").Value = "  +2.40%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.555"
$ws.Range("D22").Style = $ws.Range("C22").Style
$ws.Range("E22").Value = "  +5.54%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "11.25"
$ws.Range("D23").Style = $ws.Range("C23").Style
$ws.Range("E23").Value = "  +2.34%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.102"
$ws.Range("D24").Style = $ws.Range("C24").Style
$ws.Range("E24").Value = "  -0.11%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "158.19"
$ws.Range("D25").Style = $ws.Range("C25").Style
$ws.Range("E25").Value = "  +1.29%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "19.67"
$ws.Range("D26").Style = $ws.Range("C26").Style
$ws.Range("E26").Value = "  +1.74%  "
$ws.Range("E27").Value = "  +1.34%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "120.44"
$ws.Range("D28").Style = $ws.Range("C28").Style
$ws.Range("E28").Value = "  +2.28%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.946"
$ws.Range("D29").Style = $ws.Range("C29").Style
$ws.Range("E29").Value = "  +2.22%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.09460"
$ws.Range("D30").Style = $ws.Range("C30").Style
$ws.Range("E30").Value = "  +1.04%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.9159"
$ws.Range("D31").Style = $ws.Range("C31").Style
$ws.Range("E31").Value = "  +1.59%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.359"
$ws.Range("D32").Style = $ws.Range("C32").Style
$ws.Range("E32").Value = "  +3.10%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.286"
$ws.Range("D33").Style = $ws.Range("C33").Style
$ws.Range("E33").Value = "  +0.28%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.221"
$ws.Range("D34").Style = $ws.Range("C34").Style
$ws.Range("E34").Value = "  -0.87%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.000003505"
$ws.Range("D35").Style = $ws.Range("C35").Style
$ws.Range("E35").Value = "  +81.19%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.05846"
$ws.Range("D36").Style = $ws.Range("C36").Style
$ws.Range("E36").Value = "  +1.80%  "
$ws.Range("E37").Value = "  +0.64%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02126"
$ws.Range("D38").Style = $ws.Range("C38").Style
$ws.Range("E38").Value = "  +2.45%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "7.907"
$ws.Range("D39").Style = $ws.Range("C39").Style
$ws.Range("E39").Value = "  +3.70%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.5783"
$ws.Range("D40").Style = $ws.Range("C40").Style
$ws.Range("E40").Value = "  +2.24%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.1826"
$ws.Range("D41").Style = $ws.Range("C41").Style
$ws.Range("E41").Value = "  +2.89%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "9.897"
$ws.Range("D42").Style = $ws.Range("C42").Style
$ws.Range("E42").Value = "  +2.16%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.808"
$ws.Range("D43").Style = $ws.Range("C43").Style
$ws.Range("E43").Value = "  +10.49%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "12.11"
$ws.Range("D44").Style = $ws.Range("C44").Style
$ws.Range("E44").Value = "  +1.11%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.5407"
$ws.Range("D45").Style = $ws.Range("C45").Style
$ws.Range("E45").Value = "  +1.30%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.208"
$ws.Range("D46").Style = $ws.Range("C46").Style
$ws.Range("E46").Value = "  -1.13%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.06978"
$ws.Range("D47").Style = $ws.Range("C47").Style
$ws.Range("E47").Value = "  -0.86%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.872"
$ws.Range("D48").Style = $ws.Range("C48").Style
$ws.Range("E48").Value = "  +1.50%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "114.26"
$ws.Range("D49").Style = $ws.Range("C49").Style
$ws.Range("E49").Value = "  +1.66%  "
$ws.Range("E50").Value = "  +8.54%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "73.66"
$ws.Range("D51").Style = $ws.Range("C51").Style
$ws.Range("E51").Value = "  +4.00%  "
